# Update countries & provincias Spain
# - Swap Ucrania/Polonia rows (row 38 becomes Polonia w/ new stats,
#   row 39 becomes Ucrania w/ the stats that used to be on row 38)
# - Refresh several country stat rows (48, 55, 98) with newer numbers
# - Bump the "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 10:35"

# Row 38 -> Polonia (new figures)
$ws.Range("A38").Value = "Polonia"
$ws.Range("B38").Value = 23376
$ws.Range("C38").Value = 221
$ws.Range("D38").Value = 11016
$ws.Range("E38").Value = 11309
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1051

# Row 39 -> Ucrania (old row-38 figures)
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 23204
$ws.Range("C39").Value = 393
$ws.Range("D39").Value = 9311
$ws.Range("E39").Value = 13197
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 17
$ws.Range("H39").Value = 696

# Row 48 refreshed figures
$ws.Range("B48").Value = 14525
$ws.Range("C48").Value = 866
$ws.Range("E48").Value = 12973
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 249

# Row 55 refreshed figures
$ws.Range("E55").Value = 7382
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 42

# Row 98 refreshed figures
$ws.Range("B98").Value = 1670
$ws.Range("C98").Value = 8
$ws.Range("D98").Value = 1229
$ws.Range("E98").Value = 371
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 70
